# Natmi following Dr Hou advice
# Update LR-pairs data for Robo2-Robo2 (natmi output) to include the new "ECs" cluster,
# expanding the Sending/Target cluster grid from 2x2 (FAPs/sCs) to a full 3x3 (ECs/FAPs/sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Each inner array is one data row (columns A..T) for sheet rows 2..10,
# covering every Sending-cluster x Target-cluster combination of ECs/FAPs/sCs
# for the Robo2 ligand/receptor pair.
$rowData = @(
    @("ECs", "Robo2", "Robo2", "ECs", 2, 0.6666666666666666, 1.054425, 3.163275, 0.5590543736060202, 0.5590543736060202, 2, 0.6666666666666666, 1.054425, 3.163275, 0.5590543736060202, 0.5590543736060202, 1.111812080625, 10.006308725625, 0.3125417926480196, 0.3125417926480196),
    @("ECs", "Robo2", "Robo2", "FAPs", 2, 0.6666666666666666, 1.054425, 3.163275, 0.5590543736060202, 0.5590543736060202, 3, 1, 0.7756663333333332, 2.326999, 0.4112569941996302, 0.4112569941996303, 0.8178819735249998, 7.360937761724998, 0.229915021283369, 0.229915021283369),
    @("ECs", "Robo2", "Robo2", "sCs", 2, 0.6666666666666666, 1.054425, 3.163275, 0.5590543736060202, 0.5590543736060202, 3, 1, 0.05599533333333334, 0.167986, 0.0296886321943495, 0.02968863219434951, 0.05904287935000001, 0.53138591415, 0.01659755967463159, 0.01659755967463159),
    @("FAPs", "Robo2", "Robo2", "ECs", 3, 1, 0.7756663333333332, 2.326999, 0.4112569941996302, 0.4112569941996303, 2, 0.6666666666666666, 1.054425, 3.163275, 0.5590543736060202, 0.5590543736060202, 0.8178819735249998, 7.360937761724998, 0.229915021283369, 0.229915021283369),
    @("FAPs", "Robo2", "Robo2", "FAPs", 3, 1, 0.7756663333333332, 2.326999, 0.4112569941996302, 0.4112569941996303, 3, 1, 0.7756663333333332, 2.326999, 0.4112569941996302, 0.4112569941996303, 0.6016582606667776, 5.414924346000999, 0.1691323152781147, 0.1691323152781148),
    @("FAPs", "Robo2", "Robo2", "sCs", 3, 1, 0.7756663333333332, 2.326999, 0.4112569941996302, 0.4112569941996303, 3, 1, 0.05599533333333334, 0.167986, 0.0296886321943495, 0.02968863219434951, 0.04343369489044444, 0.390903254014, 0.01220965763814655, 0.01220965763814655),
    @("sCs", "Robo2", "Robo2", "ECs", 3, 1, 0.05599533333333334, 0.167986, 0.0296886321943495, 0.02968863219434951, 2, 0.6666666666666666, 1.054425, 3.163275, 0.5590543736060202, 0.5590543736060202, 0.05904287935000001, 0.53138591415, 0.01659755967463159, 0.01659755967463159),
    @("sCs", "Robo2", "Robo2", "FAPs", 3, 1, 0.05599533333333334, 0.167986, 0.0296886321943495, 0.02968863219434951, 3, 1, 0.7756663333333332, 2.326999, 0.4112569941996302, 0.4112569941996303, 0.04343369489044444, 0.390903254014, 0.01220965763814655, 0.01220965763814655),
    @("sCs", "Robo2", "Robo2", "sCs", 3, 1, 0.05599533333333334, 0.167986, 0.0296886321943495, 0.02968863219434951, 3, 1, 0.05599533333333334, 0.167986, 0.0296886321943495, 0.02968863219434951, 0.003135477355111112, 0.02821929619600001, 0.0008814148815713659, 0.0008814148815713661)
)

for ($r = 0; $r -lt $rowData.Length; $r++) {
    $rowValues = $rowData[$r]
    $sheetRow = 2 + $r
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Range($cols[$c] + $sheetRow).Value = $rowValues[$c]
    }
}
